$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4
$ws.Range("G4").Value = 1.45
$ws.Range("I4").Value = 6.5
$ws.Range("K4").Value = 2.4
$ws.Range("L4").Value = 6.5
$ws.Range("Y4").Value = 1.95
$ws.Range("Z4").Value = 1.8

# Row 5
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2.35
$ws.Range("R5").Value = 1.57
$ws.Range("S5").Value = 3.55
$ws.Range("T5").Value = 1.3
$ws.Range("U5").Value = 4.33
$ws.Range("V5").Value = 1.2

# Row 6
$ws.Range("I6").Value = 3.25
$ws.Range("J6").Value = 3.1
$ws.Range("N6").Value = 8.5
$ws.Range("AE6").Value = 21

# Row 7
$ws.Range("G7").Value = 2.3
$ws.Range("I7").Value = 3.6
$ws.Range("J7").Value = 3.2
$ws.Range("AC7").Value = 11
$ws.Range("AM7").Value = 15

# Row 8
$ws.Range("G8").Value = 2.2
$ws.Range("H8").Value = 2.75
$ws.Range("I8").Value = 4.2
$ws.Range("J8").Value = 3.2
$ws.Range("K8").Value = 1.73
$ws.Range("L8").Value = 5
$ws.Range("M8").Value = 1.18
$ws.Range("N8").Value = 4.5
$ws.Range("O8").Value = 1.8
$ws.Range("P8").Value = 1.91
$ws.Range("Q8").Value = 3.6
$ws.Range("R8").Value = 1.29
$ws.Range("S8").Value = 7.2
$ws.Range("T8").Value = 1.1
$ws.Range("U8").Value = 9
$ws.Range("V8").Value = 1.07
$ws.Range("W8").Value = 1.83
$ws.Range("X8").Value = 1.98
$ws.Range("Y8").Value = 2.75
$ws.Range("Z8").Value = 1.4
$ws.Range("AA8").Value = 4.5
$ws.Range("AB8").Value = 8.5
$ws.Range("AC8").Value = 12
$ws.Range("AD8").Value = 21
$ws.Range("AE8").Value = 29
$ws.Range("AG8").Value = 4.33
$ws.Range("AI8").Value = 29
$ws.Range("AJ8").Value = 151
$ws.Range("AL8").Value = 7
$ws.Range("AP8").Value = 51
$ws.Range("AR8").Value = 2.8
$ws.Range("AS8").Value = 1.44

# Row 9
$ws.Range("O9").Value = 1.5
$ws.Range("P9").Value = 2.5
$ws.Range("Q9").Value = 2.5
$ws.Range("R9").Value = 1.5
$ws.Range("S9").Value = 4.1
$ws.Range("T9").Value = 1.24
$ws.Range("AR9").Value = 1.9
$ws.Range("AS9").Value = 1.95

# Row 26
$ws.Range("G26").Value = 6.5
$ws.Range("H26").Value = 4.33
$ws.Range("I26").Value = 1.45
$ws.Range("K26").Value = 2.3
$ws.Range("L26").Value = 2
$ws.Range("M26").Value = 1.05
$ws.Range("N26").Value = 11
$ws.Range("O26").Value = 1.25
$ws.Range("P26").Value = 3.75
$ws.Range("Q26").Value = 1.85
$ws.Range("R26").Value = 2
$ws.Range("U26").Value = 3
$ws.Range("V26").Value = 1.36
$ws.Range("W26").Value = 1.36
$ws.Range("X26").Value = 3
$ws.Range("Y26").Value = 2
$ws.Range("Z26").Value = 1.75
$ws.Range("AB26").Value = 34
$ws.Range("AG26").Value = 11
$ws.Range("AK26").Value = 401
$ws.Range("AQ26").Value = 29
